$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("report")

# "Task Code" -> "Job Code" (column header, cell B3)
$ws.Range("B3").Value = "Job Code"

# "Economic Activity Index update" -> "Economic Activity update" (cell D4)
$ws.Range("D4").Value = "Economic Activity update"

# "aggregate Economic Activities: has a format error on energia,água e san"
#   -> "Could not fetch the Economic Activities Index, the url is http://www.ine.gov.mz/estatisticas/estatisticas-economicas/icce" (cell F4)
$ws.Range("F4").Value = "Could not fetch the Economic Activities Index, the url is http://www.ine.gov.mz/estatisticas/estatisticas-economicas/icce"

# Update timestamp value in G4
$ws.Range("G4").Value = 44831.64273728427
